$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.900.69'
$ws.Range('E2').Value = '  -0.13%  '
$ws.Range('D3').Value = '1.639.65'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.011'
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.06'
$ws.Range('E5').Value = '  -0.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5035'
$ws.Range('E6').Value = '  +0.10%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.007'
$ws.Range('E7').Value = '  -0.75%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2566'
$ws.Range('E8').Value = '  -0.37%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06367'
$ws.Range('E9').Value = '  -0.50%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.44'
$ws.Range('E10').Value = '  -0.18%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07799'
$ws.Range('E11').Value = '  +0.45%  '
$ws.Range('D12').Value = '1.649.64'
$ws.Range('E12').Value = '  +1.08%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.250'
$ws.Range('E13').Value = '  -0.29%  '
$ws.Range('D14').Value = '1.864.46'
$ws.Range('E14').Value = '  +0.27%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5416'
$ws.Range('E15').Value = '  -0.34%  '
$ws.Range('D16').Value = '0.0₅7855'
$ws.Range('E16').Value = '  -1.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.43'
$ws.Range('E17').Value = '  +1.51%  '
$ws.Range('D18').Value = '25.946.37'
$ws.Range('E18').Value = '  +0.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.010'
$ws.Range('E19').Value = '  -0.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '195.13'
$ws.Range('E20').Value = '  -4.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.359'
$ws.Range('E21').Value = '  +1.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.871'
$ws.Range('E22').Value = '  -1.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.957'
$ws.Range('E23').Value = '  -0.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.007'
$ws.Range('E24').Value = '  -0.94%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.902'
$ws.Range('E25').Value = '  -3.63%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '140.27'
$ws.Range('E26').Value = '  -1.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1128'
$ws.Range('E27').Value = '  -2.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.792'
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.61'
$ws.Range('E29').Value = '  -0.59%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.241'
$ws.Range('E30').Value = '  +0.11%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.04848'
$ws.Range('E31').Value = '  -2.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.237'
$ws.Range('E32').Value = '  -0.70%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.172'
$ws.Range('E33').Value = '  -0.64%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.527'
$ws.Range('E34').Value = '  -0.53%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.373'
$ws.Range('E35').Value = '  +1.01%  '
$ws.Range('B36').Value = 'MXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.605'
$ws.Range('E36').Value = '  -0.58%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.8848'
$ws.Range('E37').Value = '  -0.53%  '
$ws.Range('D38').Value = '1.125.91'
$ws.Range('E38').Value = '  +1.22%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5498'
$ws.Range('E39').Value = '  -2.63%  '
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.006'
$ws.Range('E41').Value = '  -1.10%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.658'
$ws.Range('E42').Value = '  +0.79%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8111'
$ws.Range('E43').Value = '  -0.55%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.37'
$ws.Range('E44').Value = '  -0.37%  '
$ws.Range('D45').Value = '1.775.11'
$ws.Range('E45').Value = '  +0.20%  '
$ws.Range('E46').Value = '  +4.15%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4549'
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.003'
$ws.Range('E48').Value = '  -1.37%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '54.87'
$ws.Range('E49').Value = '  +0.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05051'
$ws.Range('E50').Value = '  +0.33%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.009'
$ws.Range('E51').Value = '  -0.51%  '
